$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 12.48634390574333
$ws.Range("R2").Value = 112.37709515169
$ws.Range("S2").Value = 0.0004587857516120906
$ws.Range("T2").Value = 0.0004587857516120906
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 31.012417031995
$ws.Range("R3").Value = 279.111753287955
$ws.Range("S3").Value = 0.00113948928243014
$ws.Range("T3").Value = 0.00113948928243014
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 4.621767122873889
$ws.Range("R4").Value = 41.595904105865
$ws.Range("S4").Value = 0.0001698175958671479
$ws.Range("T4").Value = 0.0001698175958671479
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 6733.042869632431
$ws.Range("R5").Value = 60597.38582669188
$ws.Range("S5").Value = 0.2473922035865027
$ws.Range("T5").Value = 0.2473922035865027
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.6144496936817382
$ws.Range("T6").Value = 0.6144496936817382
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.09157117260445614
$ws.Range("T7").Value = 0.09157117260445614
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 316.5124999569534
$ws.Range("R8").Value = 2848.61249961258
$ws.Range("S8").Value = 0.01162961922909875
$ws.Range("T8").Value = 0.01162961922909875
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("Q9").Value = 786.12424250059
$ws.Range("R9").Value = 7075.118182505309
$ws.Range("S9").Value = 0.02888456414292937
$ws.Range("T9").Value = 0.02888456414292937
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.004304654125365513
$ws.Range("T10").Value = 0.004304654125365513
